$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 2 4 '27.052.72'
Set-TextValue 2 5 '  -2.02%  '
Set-TextValue 3 4 '1.798.01'
Set-TextValue 3 5 '  -2.63%  '
Set-TextValue 4 5 '  -0.47%  '
Set-TextValue 5 4 '307.51'
Set-TextValue 5 5 '  -2.72%  '
Set-TextValue 6 5 '  -0.28%  '
Set-TextValue 7 4 '0.4202'
Set-TextValue 7 5 '  -2.46%  '
Set-TextValue 8 5 '  -3.05%  '
Set-TextValue 9 4 '0.07102'
Set-TextValue 9 5 '  -3.46%  '
Set-TextValue 10 4 '0.8437'
Set-TextValue 10 5 '  -3.84%  '
Set-TextValue 11 4 '20.15'
Set-TextValue 11 5 '  -4.31%  '
Set-TextValue 12 4 '1.788.29'
Set-TextValue 12 5 '  -5.80%  '
Set-TextValue 13 4 '5.289'
Set-TextValue 13 5 '  -3.66%  '
Set-TextValue 14 4 '6.366'
Set-TextValue 14 5 '  -3.70%  '
Set-TextValue 15 4 '0.06759'
Set-TextValue 15 5 '  -2.91%  '
Set-TextValue 16 4 '1.002'
Set-TextValue 16 5 '  -0.71%  '
Set-TextValue 17 4 '80.19'
Set-TextValue 17 5 '  -1.81%  '
Set-TextValue 18 4 '0.000008667'
Set-TextValue 18 5 '  -4.41%  '
Set-TextValue 19 5 '  -0.09%  '
Set-TextValue 20 5 '  -3.65%  '
Set-TextValue 21 4 '27.060.94'
Set-TextValue 21 5 '  -2.78%  '
Set-TextValue 22 4 '5.053'
Set-TextValue 22 5 '  -0.60%  '
Set-TextValue 23 4 '10.99'
Set-TextValue 23 5 '  +0.08%  '
Set-TextValue 24 4 '2.020.52'
Set-TextValue 24 5 '  -3.82%  '
Set-TextValue 25 4 '1.922'
Set-TextValue 25 5 '  -3.45%  '
Set-TextValue 26 4 '152.73'
Set-TextValue 26 5 '  -1.23%  '
Set-TextValue 27 4 '18.12'
Set-TextValue 27 5 '  -5.21%  '
Set-TextValue 28 4 '5.014'
Set-TextValue 28 5 '  -6.23%  '
Set-TextValue 29 4 '113.00'
Set-TextValue 29 5 '  -2.69%  '
Set-TextValue 30 4 '1.642'
Set-TextValue 30 5 '  -12.43%  '
Set-TextValue 31 4 '0.08997'
Set-TextValue 31 5 '  +0.72%  '
Set-TextValue 32 5 '  -8.24%  '
Set-TextValue 33 4 '2.864'
Set-TextValue 33 5 '  -4.00%  '
Set-TextValue 34 4 '4.328'
Set-TextValue 34 5 '  -6.36%  '
Set-TextValue 35 4 '1.085'
Set-TextValue 35 5 '  -7.93%  '
Set-TextValue 37 5 '  -3.01%  '
Set-TextValue 38 4 '0.05127'
Set-TextValue 38 5 '  -5.96%  '
Set-TextValue 39 4 '0.01900'
Set-TextValue 39 5 '  -3.12%  '
Set-TextValue 40 4 '0.1626'
Set-TextValue 40 5 '  -4.01%  '
Set-TextValue 41 4 '0.4958'
Set-TextValue 41 5 '  -4.49%  '
Set-TextValue 42 4 '2.605'
Set-TextValue 42 5 '  -8.08%  '
Set-TextValue 43 4 '8.029'
Set-TextValue 43 5 '  -7.17%  '
Set-TextValue 44 4 '5.906'
Set-TextValue 44 5 '  -12.78%  '
Set-TextValue 45 4 '104.95'
Set-TextValue 45 5 '  -1.83%  '
Set-TextValue 46 2 'PaxDollar'
Set-TextValue 46 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 46 4 '1.001'
Set-TextValue 46 5 '  -0.08%  '
Set-TextValue 47 2 'EnergySwap'
Set-TextValue 47 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 47 4 '10.15'
Set-TextValue 47 5 '  -4.96%  '
Set-TextValue 48 4 '0.06292'
Set-TextValue 48 5 '  -4.11%  '
Set-TextValue 49 4 '0.4525'
Set-TextValue 49 5 '  -5.86%  '
Set-TextValue 50 4 '1.599'
Set-TextValue 50 5 '  -4.30%  '
Set-TextValue 51 5 '  -8.49%  '
